$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("A7").Value = 220622
$ws.Range("F7").Value = 0.89100000000000001
$ws.Range("G7").Value = 0.90882354974746704
$ws.Range("H7").Value = 0.36499999999999999
$ws.Range("I7").Value = 0.32725152373313898

# Row 8
$ws.Range("A8").Value = 220622
$ws.Range("F8").Value = 0.83440000000000003
$ws.Range("G8").Value = 0.86307191848754805
$ws.Range("H8").Value = 0.55679000000000001
$ws.Range("I8").Value = 0.48514580726623502

# Row 9
$ws.Range("A9").Value = 220622
$ws.Range("F9").Value = 0.86370000000000002
$ws.Range("G9").Value = 0.87181371450424106
$ws.Range("H9").Value = 1.7292000000000001
$ws.Range("I9").Value = 0.811184883117675

# Update the active selection to E8, matching the saved view state.
$ws.Range("E8").Select()
